$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "13"

$data = @(
    @("Judith Fernández Franco", "Green Bay - Packers", "Kansas City - Chiefs", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "Atlanta - Falcons", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Seattle - Seahawks (Local)", "Buffalo - Bills", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 52),
    @("Guiorguio Guitto Sánchez Muñoz", "Detroit - Lions (Local)", "Dallas - Cowboys (Local)", "Cincinnati - Bengals", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "New York - Jets (Local)", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Seattle - Seahawks (Local)", "Pittsburgh - Steelers (Local)", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 52),
    @("Neima Dorely Meza Jáquez", "Green Bay - Packers", "Kansas City - Chiefs", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "Cleveland - Browns (Local)", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "New York - Jets (Local)", "Tampa Bay - Buccaneers (Local)", "Carolina - Panthers (Local)", "Seattle - Seahawks (Local)", "Buffalo - Bills", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 41),
    @("Kenia Estefanía Cruz García", "Green Bay - Packers", "Dallas - Cowboys (Local)", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Houston - Texans", "Miami - Dolphins (Local)", "Atlanta - Falcons", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Minnesota - Vikings", "Pittsburgh - Steelers (Local)", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 49),
    @("Jaime Antonio García Morales", "Detroit - Lions (Local)", "Dallas - Cowboys (Local)", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "Atlanta - Falcons", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Seattle - Seahawks (Local)", "Buffalo - Bills", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 35),
    @("Arelí Olivas Alcántar", "Detroit - Lions (Local)", "Kansas City - Chiefs", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "Atlanta - Falcons", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Seattle - Seahawks (Local)", "Buffalo - Bills", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 53),
    @("Alberto Vargas Ortiz", "Detroit - Lions (Local)", "Kansas City - Chiefs", "Baltimore - Ravens (Local)", "Philadelphia - Eagles (Local)", "San Francisco - 49ers", "Jacksonville - Jaguars", "Indianapolis - Colts (Local)", "Miami - Dolphins (Local)", "Atlanta - Falcons", "Tampa Bay - Buccaneers (Local)", "Los Angeles - Rams", "Seattle - Seahawks (Local)", "Buffalo - Bills", "Los Angeles - Chargers (Local)", "Denver - Broncos", "New England - Patriots (Local)", 47)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt $data[$r].Count; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

Write-Host "Done. Sheets:" $wb.Worksheets.Count